$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.932.43'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '3.394.87'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.11'
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.76'
$ws.Range('E6').Value = '  -1.69%  '
$ws.Range('D7').Value = '3.395.00'
$ws.Range('E7').Value = '  -1.47%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.56'
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('E12').Value = '  +1.36%  '
$ws.Range('D13').Value = '3.972.96'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.18'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('D17').Value = '3.400.04'
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').Value = '61.005.19'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.85'
$ws.Range('E20').Value = '  -3.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.97'
$ws.Range('E21').Value = '  -5.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '383.10'
$ws.Range('E22').Value = '  -4.92%  '
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.26'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('E26').Value = '  -4.34%  '
$ws.Range('D27').Value = '3.530.17'
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.40'
$ws.Range('E30').Value = '  -2.95%  '
$ws.Range('E31').Value = '  -3.57%  '
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.42'
$ws.Range('E33').Value = '  -2.93%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.52'
$ws.Range('E35').Value = '  -1.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.00'
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '167.73'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('D38').Value = '3.425.08'
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.00'
$ws.Range('E39').Value = '  -3.04%  '
$ws.Range('E40').Value = '  -4.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.79'
$ws.Range('E41').Value = '  +2.48%  '
$ws.Range('E42').Value = '  -2.46%  '
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').Value = '2.477.31'
$ws.Range('E48').Value = '  -5.06%  '
$ws.Range('E49').Value = '  -1.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.99'
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('E51').Value = '  +1.14%  '
